$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look like plain numbers,
# so Excel stores them as text (matching the source data) instead of
# auto-converting them to numeric values.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values row by row.
# Row 2
$ws.Range("D2").Value = '29.739.65'
$ws.Range("E2").Value = '  +2.86%  '
# Row 3
$ws.Range("D3").Value = '1.864.96'
$ws.Range("E3").Value = '  +2.04%  '
# Row 4
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.00%  '
# Row 5
$ws.Range("D5").Value = '246.53'
$ws.Range("E5").Value = '  +2.80%  '
# Row 6
$ws.Range("D6").Value = '0.7018'
$ws.Range("E6").Value = '  +2.34%  '
# Row 7
$ws.Range("E7").Value = '  -0.01%  '
# Row 8
$ws.Range("D8").Value = '0.07777'
$ws.Range("E8").Value = '  +2.10%  '
# Row 9
$ws.Range("D9").Value = '0.3084'
$ws.Range("E9").Value = '  +2.30%  '
# Row 10
$ws.Range("D10").Value = '23.80'
$ws.Range("E10").Value = '  +1.29%  '
# Row 11
$ws.Range("D11").Value = '0.07843'
$ws.Range("E11").Value = '  +1.18%  '
# Row 12
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.875.34'
$ws.Range("E12").Value = '  +2.53%  '
# Row 13
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '5.183'
$ws.Range("E13").Value = '  +2.69%  '
# Row 14
$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D14").Value = '92.97'
$ws.Range("E14").Value = '  +3.03%  '
# Row 15
$ws.Range("D15").Value = '0.6962'
$ws.Range("E15").Value = '  +3.47%  '
# Row 16
$ws.Range("D16").Value = '6.646'
$ws.Range("E16").Value = '  +3.27%  '
# Row 17
$ws.Range("D17").Value = '29.741.65'
$ws.Range("E17").Value = '  +2.87%  '
# Row 18
$ws.Range("D18").Value = '0.000008397'
$ws.Range("E18").Value = '  +1.59%  '
# Row 19
$ws.Range("D19").Value = '2.117.97'
$ws.Range("E19").Value = '  +1.59%  '
# Row 20
$ws.Range("D20").Value = '244.29'
$ws.Range("E20").Value = '  +0.45%  '
# Row 21
$ws.Range("E21").Value = '  +1.59%  '
# Row 22
$ws.Range("E22").Value = '  +0.01%  '
# Row 23
$ws.Range("D23").Value = '7.665'
$ws.Range("E23").Value = '  +3.58%  '
# Row 24
$ws.Range("E24").Value = '  -0.03%  '
# Row 25
$ws.Range("D25").Value = '0.1521'
$ws.Range("E25").Value = '  +3.57%  '
# Row 26
$ws.Range("D26").Value = '8.979'
$ws.Range("E26").Value = '  +3.02%  '
# Row 27
$ws.Range("D27").Value = '160.49'
$ws.Range("E27").Value = '  -0.35%  '
# Row 28
$ws.Range("E28").Value = '  +1.51%  '
# Row 29
$ws.Range("E29").Value = '  +1.21%  '
# Row 30
$ws.Range("D30").Value = '4.284'
$ws.Range("E30").Value = '  +1.96%  '
# Row 31
$ws.Range("D31").Value = '4.211'
$ws.Range("E31").Value = '  +1.31%  '
# Row 32
$ws.Range("D32").Value = '1.199'
$ws.Range("E32").Value = '  +0.85%  '
# Row 33
$ws.Range("E33").Value = '  +0.11%  '
# Row 34
$ws.Range("E34").Value = '  +4.46%  '
# Row 35
$ws.Range("D35").Value = '1.924'
$ws.Range("E35").Value = '  +6.04%  '
# Row 36
$ws.Range("E36").Value = '  +1.88%  '
# Row 37
$ws.Range("D37").Value = '2.696'
$ws.Range("E37").Value = '  -0.06%  '
# Row 38
$ws.Range("D38").Value = '1.345.53'
$ws.Range("E38").Value = '  +10.65%  '
# Row 39
$ws.Range("D39").Value = '0.01890'
$ws.Range("E39").Value = '  +3.18%  '
# Row 40
$ws.Range("D40").Value = '2.742'
$ws.Range("E40").Value = '  +1.96%  '
# Row 41
$ws.Range("D41").Value = '0.9645'
$ws.Range("E41").Value = '  +6.02%  '
# Row 42
$ws.Range("E42").Value = '  +13.22%  '
# Row 43
$ws.Range("D43").Value = '106.76'
$ws.Range("E43").Value = '  -2.03%  '
# Row 44
$ws.Range("E44").Value = '  +0.04%  '
# Row 45
$ws.Range("D45").Value = '9.797'
$ws.Range("E45").Value = '  +3.51%  '
# Row 46
$ws.Range("D46").Value = '2.012.15'
$ws.Range("E46").Value = '  +0.94%  '
# Row 47
$ws.Range("E47").Value = '  +2.30%  '
# Row 48
$ws.Range("D48").Value = '65.57'
$ws.Range("E48").Value = '  +4.31%  '
# Row 49
$ws.Range("E49").Value = '  +0.57%  '
# Row 50
$ws.Range("D50").Value = '1.794'
$ws.Range("E50").Value = '  +4.07%  '
# Row 51
$ws.Range("D51").Value = '7.040'
$ws.Range("E51").Value = '  +2.38%  '
